# Zeitaufzeichnung update: rework February work-log entries
# (split combined time ranges into individual rows, add two new
# entries for slide work / ansible_facts exercise, update totals).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Summary row -------------------------------------------------------
# (A2's label text itself doesn't change - only the SUM range below does)

# --- Make room: insert two fresh rows before the current row 11 ------
# (old rows 11.. shift down to 13.., giving us 11/12 free for new
# entries and 15/16 free once the old placeholder rows shift down)
$ws.Rows("11:12").Insert()

# Set the summary formula after the insert so the range literal itself
# reads A5:A60 (it would auto-grow to A5:A62 if set beforehand).
$ws.Range("C2").Formula = "=SUM(A5:A60)"

# --- Row 5 - 10: only the Zeitraum text actually changes; Tasks (D) and
#     row 6's Zeitraum (C6) are untouched by this commit, so leave those
#     cells alone (avoids needlessly re-encoding unrelated shared
#     strings, e.g. their CRLF line breaks). ----------------------------
$ws.Range("C5").Value = "20:00 - 23:0"

$ws.Range("A7").Value = 0.5
$ws.Range("B7").Value = 43491
$ws.Range("C7").Value = "15:30 - 16:00"
$ws.Range("D7").Value = "Übersicht"

$ws.Range("A8").Value = 1.75
$ws.Range("B8").Value = 43491
$ws.Range("C8").Value = "19:30 - 21:15"
$ws.Range("D8").Value = "Übersicht"

$ws.Range("A9").Value = 4
$ws.Range("B9").Value = 43496
$ws.Range("C9").Value = "20:00 - 24:00"
$ws.Range("D9").Value = "Präsentation"

$ws.Range("A10").Value = 4
$ws.Range("B10").Value = 43499
$ws.Range("C10").Value = "19:00 - 24:00"
$ws.Range("D10").Value = "AWS Environment"

# --- Row 11 & 12: brand-new rows (inserted above) ---------------------
$ws.Range("A11").Value = 1
$ws.Range("B11").Value = 43502
$ws.Range("C11").Value = "19:00 - 20:00"
$ws.Range("D11").Value = "AWS Environment"

$ws.Range("A12").Value = 2.5
$ws.Range("B12").Value = 43502
$ws.Range("C12").Value = "22:30 - 01:00"
$ws.Range("D12").Value = "AWS Environment"

# --- Row 13 & 14: former rows 11 & 12 (shifted down) -------------------
$ws.Range("A13").Value = 1.5
$ws.Range("B13").Value = 43503
$ws.Range("C13").Value = "16:30 - 18:00"
$ws.Range("D13").Value = "AWS Env - Config AnsibleHost/Smoketest"

$ws.Range("A14").Value = 3
$ws.Range("B14").Value = 43503
$ws.Range("C14").Value = "20:30 - 23:30"
$ws.Range("D14").Value = "WinSCP Syncing, vscode Task Sequences"

# --- Row 15 & 16: former blank placeholder rows, now used --------------
$ws.Range("B15").NumberFormat = "m/d/yy"
$ws.Range("C15").NumberFormat = "h:mm"
$ws.Range("A15").Value = 1
$ws.Range("B15").Value = 43504
$ws.Range("C15").Value = "20:15 - 21:15"
$ws.Range("D15").Value = "Präsentation"

$ws.Range("B16").NumberFormat = "m/d/yy"
$ws.Range("C16").NumberFormat = "h:mm"
$ws.Range("A16").Value = 4.5
$ws.Range("B16").Value = 43507
$ws.Range("C16").Value = "19:30 - 00:00"
$ws.Range("D16").Value = "Präsentation/Übungs 1 - ansible_facts/aws fixes"

# --- Selection, matching the saved cursor position in the workbook ----
$ws.Range("D17").Select()
